# Commit: "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The worker/period table (rows 16-88) is re-sorted: instead of being
# grouped by worker (each worker's 7-8 periods together), the data is
# now grouped by period (Periodo Mora, column E) in ascending order
# (1910 -> 2005), and within each period the workers appear in a fixed
# order. Column G (Salario Basico) is normalized to 828116 for every
# row. Only the cells whose value actually changes are touched below;
# column B ("CC") is identical in both versions so it is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1121334531"
$ws.Range("D16").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("E16").Value = "1910"
$ws.Range("F16").Value = 33125
$ws.Range("C17").Value = "1052216007"
$ws.Range("D17").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E17").Value = "1910"
$ws.Range("C18").Value = "92188076"
$ws.Range("D18").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("E18").Value = "1910"
$ws.Range("C19").Value = "5047742"
$ws.Range("D19").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E19").Value = "1910"
$ws.Range("E20").Value = "1911"
$ws.Range("F20").Value = 17667
$ws.Range("C21").Value = "1121334531"
$ws.Range("D21").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("E21").Value = "1911"
$ws.Range("C22").Value = "5091810"
$ws.Range("D22").Value = "JOSE MIGUEL ESCOBAR FLORIAN"
$ws.Range("C23").Value = "1052216007"
$ws.Range("D23").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E23").Value = "1911"
$ws.Range("F23").Value = 33125
$ws.Range("G23").Value = 828116
$ws.Range("C24").Value = "92188076"
$ws.Range("D24").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("E24").Value = "1911"
$ws.Range("G24").Value = 828116
$ws.Range("C25").Value = "5047742"
$ws.Range("D25").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E25").Value = "1911"
$ws.Range("G25").Value = 828116
$ws.Range("C26").Value = "1052571154"
$ws.Range("D26").Value = "JAIME LUIS CAMPO CASTRO"
$ws.Range("E26").Value = "1911"
$ws.Range("F26").Value = 17667
$ws.Range("G26").Value = 828116
$ws.Range("C27").Value = "1049897982"
$ws.Range("D27").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("E27").Value = "1911"
$ws.Range("F27").Value = 17667
$ws.Range("G27").Value = 828116
$ws.Range("C28").Value = "12584600"
$ws.Range("D28").Value = "MAURICIO CUADROS TORRES"
$ws.Range("E28").Value = "1911"
$ws.Range("F28").Value = 17667
$ws.Range("G28").Value = 828116
$ws.Range("C29").Value = "3820419"
$ws.Range("D29").Value = "JUSTO PASTOR PALLARES MURILLO"
$ws.Range("E29").Value = "1912"
$ws.Range("G29").Value = 828116
$ws.Range("E30").Value = "1912"
$ws.Range("G30").Value = 828116
$ws.Range("E31").Value = "1912"
$ws.Range("F31").Value = 33125
$ws.Range("C32").Value = "1052216007"
$ws.Range("D32").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E32").Value = "1912"
$ws.Range("C33").Value = "92188076"
$ws.Range("D33").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("E33").Value = "1912"
$ws.Range("C34").Value = "5047742"
$ws.Range("D34").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E34").Value = "1912"
$ws.Range("C35").Value = "1052571154"
$ws.Range("D35").Value = "JAIME LUIS CAMPO CASTRO"
$ws.Range("E35").Value = "1912"
$ws.Range("C36").Value = "1049897982"
$ws.Range("D36").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("C37").Value = "12584600"
$ws.Range("D37").Value = "MAURICIO CUADROS TORRES"
$ws.Range("E37").Value = "1912"
$ws.Range("F37").Value = 33125
$ws.Range("C38").Value = "92097983"
$ws.Range("D38").Value = "MANUEL FRANCISCO DE LA ROSA PEREZ"
$ws.Range("E38").Value = "1912"
$ws.Range("F38").Value = 33125
$ws.Range("C39").Value = "3820419"
$ws.Range("D39").Value = "JUSTO PASTOR PALLARES MURILLO"
$ws.Range("E39").Value = "2001"
$ws.Range("C40").Value = "1121334531"
$ws.Range("D40").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("E40").Value = "2001"
$ws.Range("C41").Value = "5091810"
$ws.Range("D41").Value = "JOSE MIGUEL ESCOBAR FLORIAN"
$ws.Range("E41").Value = "2001"
$ws.Range("C43").Value = "92188076"
$ws.Range("D43").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("E43").Value = "2001"
$ws.Range("C44").Value = "5047742"
$ws.Range("D44").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E44").Value = "2001"
$ws.Range("C45").Value = "1052571154"
$ws.Range("D45").Value = "JAIME LUIS CAMPO CASTRO"
$ws.Range("E45").Value = "2001"
$ws.Range("C46").Value = "1049897982"
$ws.Range("D46").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("E46").Value = "2001"
$ws.Range("F46").Value = 33125
$ws.Range("C47").Value = "12584600"
$ws.Range("D47").Value = "MAURICIO CUADROS TORRES"
$ws.Range("E47").Value = "2001"
$ws.Range("C48").Value = "92097983"
$ws.Range("D48").Value = "MANUEL FRANCISCO DE LA ROSA PEREZ"
$ws.Range("E48").Value = "2001"
$ws.Range("C49").Value = "3820419"
$ws.Range("D49").Value = "JUSTO PASTOR PALLARES MURILLO"
$ws.Range("C50").Value = "1121334531"
$ws.Range("D50").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("E50").Value = "2002"
$ws.Range("C51").Value = "5091810"
$ws.Range("D51").Value = "JOSE MIGUEL ESCOBAR FLORIAN"
$ws.Range("E51").Value = "2002"
$ws.Range("C52").Value = "1052216007"
$ws.Range("D52").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E52").Value = "2002"
$ws.Range("E53").Value = "2002"
$ws.Range("E54").Value = "2002"
$ws.Range("F54").Value = 33125
$ws.Range("C55").Value = "1052571154"
$ws.Range("D55").Value = "JAIME LUIS CAMPO CASTRO"
$ws.Range("E55").Value = "2002"
$ws.Range("C56").Value = "1049897982"
$ws.Range("D56").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("E56").Value = "2002"
$ws.Range("C57").Value = "12584600"
$ws.Range("D57").Value = "MAURICIO CUADROS TORRES"
$ws.Range("C58").Value = "92097983"
$ws.Range("D58").Value = "MANUEL FRANCISCO DE LA ROSA PEREZ"
$ws.Range("E58").Value = "2002"
$ws.Range("C59").Value = "3820419"
$ws.Range("D59").Value = "JUSTO PASTOR PALLARES MURILLO"
$ws.Range("E59").Value = "2003"
$ws.Range("C60").Value = "1121334531"
$ws.Range("D60").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("E60").Value = "2003"
$ws.Range("C61").Value = "5091810"
$ws.Range("D61").Value = "JOSE MIGUEL ESCOBAR FLORIAN"
$ws.Range("E61").Value = "2003"
$ws.Range("C62").Value = "1052216007"
$ws.Range("D62").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E62").Value = "2003"
$ws.Range("F62").Value = 33125
$ws.Range("C63").Value = "92188076"
$ws.Range("D63").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("E63").Value = "2003"
$ws.Range("C64").Value = "5047742"
$ws.Range("D64").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E65").Value = "2003"
$ws.Range("C66").Value = "1049897982"
$ws.Range("D66").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("E66").Value = "2003"
$ws.Range("C67").Value = "12584600"
$ws.Range("D67").Value = "MAURICIO CUADROS TORRES"
$ws.Range("E67").Value = "2003"
$ws.Range("C68").Value = "92097983"
$ws.Range("D68").Value = "MANUEL FRANCISCO DE LA ROSA PEREZ"
$ws.Range("E68").Value = "2003"
$ws.Range("F68").Value = 33125
$ws.Range("C69").Value = "3820419"
$ws.Range("D69").Value = "JUSTO PASTOR PALLARES MURILLO"
$ws.Range("E69").Value = "2004"
$ws.Range("F69").Value = 33125
$ws.Range("G69").Value = 828116
$ws.Range("C70").Value = "1121334531"
$ws.Range("D70").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("G70").Value = 828116
$ws.Range("C71").Value = "5091810"
$ws.Range("D71").Value = "JOSE MIGUEL ESCOBAR FLORIAN"
$ws.Range("E71").Value = "2004"
$ws.Range("G71").Value = 828116
$ws.Range("C72").Value = "1052216007"
$ws.Range("D72").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E72").Value = "2004"
$ws.Range("G72").Value = 828116
$ws.Range("C73").Value = "92188076"
$ws.Range("D73").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("E73").Value = "2004"
$ws.Range("G73").Value = 828116
$ws.Range("C74").Value = "5047742"
$ws.Range("D74").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E74").Value = "2004"
$ws.Range("G74").Value = 828116
$ws.Range("C75").Value = "1052571154"
$ws.Range("D75").Value = "JAIME LUIS CAMPO CASTRO"
$ws.Range("E75").Value = "2004"
$ws.Range("F75").Value = 33125
$ws.Range("G75").Value = 828116
$ws.Range("C76").Value = "1049897982"
$ws.Range("D76").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("E76").Value = "2004"
$ws.Range("F76").Value = 33125
$ws.Range("C78").Value = "92097983"
$ws.Range("D78").Value = "MANUEL FRANCISCO DE LA ROSA PEREZ"
$ws.Range("E78").Value = "2004"
$ws.Range("C79").Value = "3820419"
$ws.Range("D79").Value = "JUSTO PASTOR PALLARES MURILLO"
$ws.Range("E79").Value = "2005"
$ws.Range("F79").Value = 20979
$ws.Range("C80").Value = "1121334531"
$ws.Range("D80").Value = "WENDYS LORENA PLATA PEREZ"
$ws.Range("E80").Value = "2005"
$ws.Range("F80").Value = 20979
$ws.Range("C81").Value = "5091810"
$ws.Range("D81").Value = "JOSE MIGUEL ESCOBAR FLORIAN"
$ws.Range("E81").Value = "2005"
$ws.Range("F81").Value = 20979
$ws.Range("C82").Value = "1052216007"
$ws.Range("D82").Value = "MIGUEL EDUARDO GALVAN URRUTIA"
$ws.Range("E82").Value = "2005"
$ws.Range("F82").Value = 20979
$ws.Range("C83").Value = "92188076"
$ws.Range("D83").Value = "JOAQUIN MIGUEL GALVAN GARCIA"
$ws.Range("C84").Value = "5047742"
$ws.Range("D84").Value = "JOSE MARIA CORDOBA MUNOS"
$ws.Range("E84").Value = "2005"
$ws.Range("F84").Value = 20979
$ws.Range("C85").Value = "1052571154"
$ws.Range("D85").Value = "JAIME LUIS CAMPO CASTRO"
$ws.Range("E85").Value = "2005"
$ws.Range("F85").Value = 20979
$ws.Range("C86").Value = "1049897982"
$ws.Range("D86").Value = "BRAYNIS GARCIA MIRANDA"
$ws.Range("E86").Value = "2005"
$ws.Range("F86").Value = 20979
$ws.Range("C87").Value = "12584600"
$ws.Range("D87").Value = "MAURICIO CUADROS TORRES"
$ws.Range("E87").Value = "2005"
$ws.Range("F87").Value = 20979
$ws.Range("E88").Value = "2005"
$ws.Range("F88").Value = 20979
